$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Add the new data row (2020-06-05) to the bottom of the table
$ws.Range("A85").Value = 43987
$ws.Range("B85").Value = 531
$ws.Range("C85").Value = 91
$ws.Range("D85").Value = 390
$ws.Range("E85").Value = 144
$ws.Range("F85").Value = 49

# Match formatting of the row above (date style + centered numeric style)
$ws.Range("A84").Copy() | Out-Null
$ws.Range("A85").PasteSpecial(-4122) | Out-Null
$ws.Range("B84:C84").Copy() | Out-Null
$ws.Range("B85:C85").PasteSpecial(-4122) | Out-Null
$ws.Range("E84:F84").Copy() | Out-Null
$ws.Range("E85:F85").PasteSpecial(-4122) | Out-Null
$ws.Range("B84").Copy() | Out-Null
$ws.Range("D85").PasteSpecial(-4122) | Out-Null

# Grow the table so the new row is included
$table = $ws.ListObjects.Item("Condicion_Pacientes")
$table.Resize($ws.Range("A1:F85"))

# Update the view selection to match the final saved state
$ws.Range("C78").Select() | Out-Null
